$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# metsData sheet: shift the metabolite-id column (A4:A24) down by one entry
# (m2 was missing before; m8 was missing before) and append two new rows
# (25, 26) for m24 / m25 so the sheet now lists all 25 metabolites.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("metsData")

$ws.Range("A4").Value = "m2"
$ws.Range("A5").Value = "m4"
$ws.Range("A6").Value = "m5"
$ws.Range("A7").Value = "m6"
$ws.Range("A8").Value = "m7"
$ws.Range("A9").Value = "m8"
$ws.Range("A10").Value = "m9"
$ws.Range("A11").Value = "m10"
$ws.Range("A12").Value = "m11"
$ws.Range("A13").Value = "m12"
$ws.Range("A14").Value = "m13"
$ws.Range("A15").Value = "m14"
$ws.Range("A16").Value = "m15"
$ws.Range("A17").Value = "m16"
$ws.Range("A18").Value = "m17"
$ws.Range("A19").Value = "m18"
$ws.Range("A20").Value = "m19"
$ws.Range("A21").Value = "m20"
$ws.Range("A22").Value = "m21"
$ws.Range("A23").Value = "m22"
$ws.Range("A24").Value = "m23"

$ws.Range("A25").Value = "m24"
$ws.Range("B25").Value = 0.99
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 1.01

$ws.Range("A26").Value = "m25"
$ws.Range("B26").Value = 0.99
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1.01

# ---------------------------------------------------------------------------
# kinetics1 sheet: move its remembered selection to H39 (without leaving it
# as the active tab).
# ---------------------------------------------------------------------------
$wsKinetics = $wb.Worksheets.Item("kinetics1")
$wsKinetics.Activate()
$wsKinetics.Range("H39").Select()

# ---------------------------------------------------------------------------
# Finish with metsData as the active sheet/tab, selection on A5 (matches
# activeTab=11 in workbook.xml and the sheetView selection for metsData).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A5").Select()
